$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3, 4, 6, 7 - column F: remove duplicated "MEC-1B-Metrologia 1" entries, leave just "-"
$ws.Range("F3").Value = "-"
$ws.Range("F4").Value = "-"
$ws.Range("F6").Value = "-"
$ws.Range("F7").Value = "-"

# Row 18
$ws.Range("B18").Value = "-"
$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "[-, 'ELM-2NA-CAD']"

# Row 19
$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "-"
$ws.Range("E19").Value = "[-, 'ELM-2NA-CAD']"

# Row 20
$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "[-, 'ELM-2NA-CAD']"

# Row 21
$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "[-, 'ELM-2NA-CAD']"
